$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Find-ParagraphIndex($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -like "*$text*") {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1) Duplicate the empty "color" paragraph that follows "Setup instructions
#    for Windows" so a second, identical empty paragraph sits right before
#    the "Download and install MongoDB..." paragraph.
# ---------------------------------------------------------------------------
$idx = Find-ParagraphIndex("Setup instructions for Windows")
$colorPara = $d.Paragraphs.Item($idx + 1)
$insPoint = $d.Range($colorPara.Range.End, $colorPara.Range.End)
$xml = '<w:p ' + $wns + '><w:pPr><w:rPr><w:color w:val="1F4E79" w:themeColor="accent1" w:themeShade="80"/></w:rPr></w:pPr></w:p>'
$insPoint.InsertXML($xml) | Out-Null

# ---------------------------------------------------------------------------
# 2) "Download and install MongoDB 3.x from <link>" paragraph:
#    - bold + size 22pt (sz 44) on the paragraph mark and the first run
#    - bold + size 22pt on the hyperlink run as well
#    - a relocated "_GoBack" bookmark right at the start of the paragraph
# ---------------------------------------------------------------------------
$idx = Find-ParagraphIndex("Download and install MongoDB 3.x from")
$downloadPara = $d.Paragraphs.Item($idx)
$downloadPara.Range.Bold = 1
$downloadPara.Range.Font.Size = 22

$hl = $d.Hyperlinks.Item(1)
$hl.Range.Bold = 1
$hl.Range.Font.Size = 22

$bmStart = $downloadPara.Range.Start
$bmRange = $d.Range($bmStart, $bmStart)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------------
# 3) Insert two new empty "ListParagraph / jc=both" paragraphs right after
#    the download paragraph (before the pre-existing empty one).
# ---------------------------------------------------------------------------
$idx = Find-ParagraphIndex("Download and install MongoDB 3.x from")
$downloadPara = $d.Paragraphs.Item($idx)
$endPoint = $d.Range($downloadPara.Range.End, $downloadPara.Range.End)
$emptyListPara = '<w:p ' + $wns + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:jc w:val="both"/></w:pPr></w:p>'
$xml = $emptyListPara + $emptyListPara
$endPoint.InsertXML($xml) | Out-Null

# ---------------------------------------------------------------------------
# 4) Insert a new run right before "Double-c" in the Double-click paragraph.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$ok = $rng.Find.Execute("Double-c", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($ok) {
    $ip = $d.Range($rng.Start, $rng.Start)
    $ip.InsertBefore("Make sure you have downloaded and installed MongoDB and ")
}
